# Split the single run containing
#   Prompt user to enter their name. "Enter your name"
# into three runs, breaking the text right before/after the "E" of
# "Enter", i.e.:
#   1) Prompt user to enter their name. "
#   2) E
#   3) nter your name"
#
# Word's object model has no direct "split run" call, but toggling a
# direct-character-formatting property on a sub-range and then clearing
# it again forces Word to record that sub-range as its own run when the
# paragraph is serialized, which is how this kind of run split shows up
# in authored documents.

$d = $word.ActiveDocument

# Locate the target sentence; Find.Execute repositions this range onto
# the match ("Enter your name").
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Enter your name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $searchRange.Start

    # Isolate just the "E" character (first character of the match).
    $charE = $d.Range($matchStart, $matchStart + 1)

    # Toggling Bold on then off on this one-character range forces Word
    # to split the run at this character's boundaries without altering
    # the visible text or formatting.
    $charE.Bold = 1
    $charE.Bold = 0
}
